$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 735.5333000000001
$ws.Range("I43").Value = 613.1667
$ws.Range("J43").Value = 817.1111
$ws.Range("K43").Value = 613.1667
$ws.Range("L43").Value = 817.1111
$ws.Range("M43").Value = -544.1667
$ws.Range("N43").Value = -955.1111
$ws.Range("H62").Value = 1829.5
$ws.Range("J62").Value = 2300.8333
$ws.Range("L62").Value = 2300.8333
$ws.Range("N62").Value = -3548.8333
$ws.Range("H65").Value = 1829.5
$ws.Range("J65").Value = 2300.8333
$ws.Range("L65").Value = 11504.1665
$ws.Range("N65").Value = -17744.1665
$ws.Range("H80").Value = 4421.1333
$ws.Range("I80").Value = 573.2
$ws.Range("J80").Value = 6345.1
$ws.Range("K80").Value = 1719.6
$ws.Range("L80").Value = 19035.3
$ws.Range("M80").Value = -721.6000000000001
$ws.Range("N80").Value = -21031.3
$ws.Range("H83").Value = 4421.1333
$ws.Range("I83").Value = 573.2
$ws.Range("J83").Value = 6345.1
$ws.Range("K83").Value = 5158.8
$ws.Range("L83").Value = 57105.9
$ws.Range("M83").Value = -166.8000000000002
$ws.Range("N83").Value = -67089.89999999999
$ws.Range("H98").Value = 1042.7037
$ws.Range("I98").Value = 825.73914
$ws.Range("J98").Value = 2290.25
$ws.Range("K98").Value = 825.73914
$ws.Range("L98").Value = 2290.25
$ws.Range("M98").Value = 672.26086
$ws.Range("N98").Value = -5286.25
$ws.Range("H100").Value = 2529.9167
$ws.Range("I100").Value = 1101.6666
$ws.Range("J100").Value = 3006
$ws.Range("K100").Value = 1101.6666
$ws.Range("L100").Value = 3006
$ws.Range("M100").Value = -560.6666
$ws.Range("N100").Value = -4088
$ws.Range("H111").Value = 5241
$ws.Range("I111").Value = 6985.6
$ws.Range("K111").Value = 20956.8
$ws.Range("M111").Value = -17889.8
$ws.Range("H112").Value = 1367.027
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1367.027
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4101.081
$ws.Range("M112").ClearContents() | Out-Null
$ws.Range("N112").Value = -6317.081
$ws.Range("H122").Value = 1042.7037
$ws.Range("I122").Value = 825.73914
$ws.Range("J122").Value = 2290.25
$ws.Range("K122").Value = 2477.21742
$ws.Range("L122").Value = 6870.75
$ws.Range("M122").Value = -27.21741999999995
$ws.Range("N122").Value = -11770.75
$ws.Range("H137").Value = 593489.9
$ws.Range("I137").Value = 1544.2667
$ws.Range("J137").Value = 1638099.8
$ws.Range("K137").Value = 4632.800099999999
$ws.Range("L137").Value = 4914299.4
$ws.Range("M137").Value = -2082.800099999999
$ws.Range("N137").Value = -4919399.4

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8202.361000000001
$ws.Range("I32").Value = 6559.7856
$ws.Range("J32").Value = 22000
$ws.Range("K32").Value = 6559.7856
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = -6272.7856
$ws.Range("N32").Value = -22574
$ws.Range("H45").Value = 2126.5
$ws.Range("I45").Value = 1835.3334
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 1835.3334
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -1458.3334
$ws.Range("N45").Value = -3754
$ws.Range("H74").Value = 6776.778
$ws.Range("I74").Value = 3882.25
$ws.Range("J74").Value = 15046.857
$ws.Range("K74").Value = 3882.25
$ws.Range("L74").Value = 15046.857
$ws.Range("M74").Value = -3008.25
$ws.Range("N74").Value = -16794.857
$ws.Range("H77").Value = 6776.778
$ws.Range("I77").Value = 3882.25
$ws.Range("J77").Value = 15046.857
$ws.Range("K77").Value = 19411.25
$ws.Range("L77").Value = 75234.285
$ws.Range("M77").Value = -15043.25
$ws.Range("N77").Value = -83970.285

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 375
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 23
$ws.Range("N22").Value = -946
$ws.Range("H99").Value = 1713.4193
$ws.Range("I99").Value = 1545.4445
$ws.Range("K99").Value = 1545.4445
$ws.Range("M99").Value = -47.44450000000006
$ws.Range("H107").Value = 2290.1765
$ws.Range("I107").Value = 2132.1
$ws.Range("J107").Value = 2516
$ws.Range("K107").Value = 2132.1
$ws.Range("L107").Value = 2516
$ws.Range("M107").Value = -212.0999999999999
$ws.Range("N107").Value = -6356
$ws.Range("H118").Value = 57379.2
$ws.Range("J118").Value = 57379.2
$ws.Range("L118").Value = 57379.2
$ws.Range("N118").Value = -60693.2

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2519.8918
$ws.Range("I31").Value = 1829.742
$ws.Range("J31").Value = 6085.6665
$ws.Range("K31").Value = 1829.742
$ws.Range("L31").Value = 6085.6665
$ws.Range("M31").Value = -1534.742
$ws.Range("N31").Value = -6675.6665
$ws.Range("H34").Value = 2519.8918
$ws.Range("I34").Value = 1829.742
$ws.Range("J34").Value = 6085.6665
$ws.Range("K34").Value = 1829.742
$ws.Range("L34").Value = 6085.6665
$ws.Range("M34").Value = -1627.742
$ws.Range("N34").Value = -6489.6665
$ws.Range("H47").Value = 18224.125
$ws.Range("J47").Value = 18224.125
$ws.Range("L47").Value = 18224.125
$ws.Range("N47").Value = -19356.125
$ws.Range("H53").Value = 42606.5
$ws.Range("J53").Value = 42606.5
$ws.Range("L53").Value = 42606.5
$ws.Range("N53").Value = -43820.5
$ws.Range("H107").Value = 399.68
$ws.Range("I107").Value = 402.10526
$ws.Range("J107").Value = 392
$ws.Range("K107").Value = 402.10526
$ws.Range("L107").Value = 392
$ws.Range("M107").Value = 1517.89474
$ws.Range("N107").Value = -4232
$ws.Range("H111").Value = 69900
$ws.Range("J111").Value = 69900
$ws.Range("L111").Value = 69900
$ws.Range("N111").Value = -78080

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 6816
$ws.Range("J120").Value = 7182.4
$ws.Range("L120").Value = 21547.2
$ws.Range("N120").Value = -31223.2
$ws.Range("H132").Value = 1896.7368
$ws.Range("I132").Value = 2237.375
$ws.Range("J132").Value = 1649
$ws.Range("K132").Value = 20136.375
$ws.Range("L132").Value = 14841
$ws.Range("M132").Value = -17606.375
$ws.Range("N132").Value = -19901

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2211.182
$ws.Range("I113").Value = 2184.6924
$ws.Range("J113").Value = 2249.4443
$ws.Range("K113").Value = 2184.6924
$ws.Range("L113").Value = 2249.4443
$ws.Range("M113").Value = -14.69239999999991
$ws.Range("N113").Value = -6589.4443
$ws.Range("H122").Value = 4423.5293
$ws.Range("I122").Value = 5315.385
$ws.Range("J122").Value = 1525
$ws.Range("K122").Value = 15946.155
$ws.Range("L122").Value = 4575
$ws.Range("M122").Value = -13496.155
$ws.Range("N122").Value = -9475
$ws.Range("H134").Value = 38074.617
$ws.Range("J134").Value = 38074.617
$ws.Range("L134").Value = 114223.851
$ws.Range("N134").Value = -119293.851
$ws.Range("H136").Value = 7575.5884
$ws.Range("J136").Value = 7575.5884
$ws.Range("L136").Value = 22726.7652
$ws.Range("N136").Value = -27826.7652

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1236.6
$ws.Range("I16").Value = 1536.8334
$ws.Range("J16").Value = 786.25
$ws.Range("K16").Value = 1536.8334
$ws.Range("L16").Value = 786.25
$ws.Range("M16").Value = -1366.8334
$ws.Range("N16").Value = -1126.25
$ws.Range("H40").Value = 3566.7026
$ws.Range("I40").Value = 3504.32
$ws.Range("K40").Value = 3504.32
$ws.Range("M40").Value = -3368.32
$ws.Range("H82").Value = 1638.7142
$ws.Range("I82").Value = 1361.5
$ws.Range("J82").Value = 2008.3334
$ws.Range("K82").Value = 1361.5
$ws.Range("L82").Value = 2008.3334
$ws.Range("M82").Value = -1000.5
$ws.Range("N82").Value = -2730.3334
$ws.Range("H85").Value = 1638.7142
$ws.Range("I85").Value = 1361.5
$ws.Range("J85").Value = 2008.3334
$ws.Range("K85").Value = 1361.5
$ws.Range("L85").Value = 2008.3334
$ws.Range("M85").Value = -113.5
$ws.Range("N85").Value = -4504.3334
$ws.Range("H122").Value = 4429.919
$ws.Range("I122").Value = 4066.7407
$ws.Range("J122").Value = 5410.5
$ws.Range("K122").Value = 12200.2221
$ws.Range("L122").Value = 16231.5
$ws.Range("M122").Value = -9750.222099999999
$ws.Range("N122").Value = -21131.5
$ws.Range("H136").Value = 4767.7896
$ws.Range("I136").Value = 3453.8667
$ws.Range("J136").Value = 6227.7036
$ws.Range("K136").Value = 10361.6001
$ws.Range("L136").Value = 18683.1108
$ws.Range("M136").Value = -7811.6001
$ws.Range("N136").Value = -23783.1108
$ws.Range("H140").Value = 53233.816
$ws.Range("J140").Value = 53233.816
$ws.Range("L140").Value = 53233.816
$ws.Range("N140").Value = -63593.816

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1184.0416
$ws.Range("I107").Value = 498.77777
$ws.Range("J107").Value = 1595.2
$ws.Range("K107").Value = 1496.33331
$ws.Range("L107").Value = 4785.6
$ws.Range("M107").Value = 423.66669
$ws.Range("N107").Value = -8625.6
$ws.Range("H113").Value = 1148.8182
$ws.Range("I113").Value = 450.33334
$ws.Range("K113").Value = 1351.00002
$ws.Range("M113").Value = 818.9999800000001
